$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: force text format ("@") before assignment so that
# numeric-looking strings (e.g. "302.91", "1.11%") remain stored as
# literal text, matching the original inline-string cell contents.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "302.91"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.11%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "32.16"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.18%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.883"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.99%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07873"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.02%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.074"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-6.23%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.839"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.25%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.843"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.53%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.05%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.73%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07803"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5.82%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08565"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-7.42%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03162"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "4.69%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.40%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001514"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.10%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005723"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.19%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2,109.78%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.464"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.43%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.117"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-6.55%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3271"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.02%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.37%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.306"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.99%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "16.92%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04582"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.80%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001225"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.24%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004456"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.15%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "4.30%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01728"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.88%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04809"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.33%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007492"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "7.89%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1364"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.46%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002361"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "7.90%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "9.36%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006144"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.99%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.09%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-61.11%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.8234"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-28.95%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.09%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.09%"
